$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: formula gains an extra " + 3" term (6 + 4 -> 6 + 4 + 3 = 13)
$ws.Range("B16").Formula = "= 6 + 4 + 3"

# New row 17: week label + formula "= 6"
# (set content first, then copy formatting from row 16 so the dependency
# graph picks the new cells up before the format-only paste touches them)
$ws.Range("A17").Value = "week 22-28/02/2016"
$ws.Range("A16").Copy()
$ws.Range("A17").PasteSpecial(-4122)

$ws.Range("B17").Formula = "= 6"
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)

# New row 18: empty but formatted cell B18 (matches column formatting)
$ws.Range("B18").Formula = "= 6"
$ws.Range("B16").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").ClearContents()

# Selection moves to E3
$ws.Range("E3").Select()
